{"js": "// Insert three new paragraphs at the end of the document body (after the\n// trailing empty paragraph, before the section break), matching the diff:\n//   MS3\n//   Finished Postal Frog model and created alternate \"Leap\" model\n//   Acquired non-functioning first pass character controller\nconst body = context.document.body;\n\nbody.insertParagraph(\"MS3\", \"End\");\nbody.insertParagraph(\"Finished Postal Frog model and created alternate \\u201CLeap\\u201D model\", \"End\");\nbody.insertParagraph(\"Acquired non-functioning first pass character controller\", \"End\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Smart/curly quotes used around \"Leap\" in the second new paragraph.\n$lq = [char]0x201C\n$rq = [char]0x201D\n\n$newParagraphs = @(\n    \"MS3\",\n    \"Finished Postal Frog model and created alternate \" + $lq + \"Leap\" + $rq + \" model\",\n    \"Acquired non-functioning first pass character controller\"\n)\n\nforeach ($text in $newParagraphs) {\n    $r = $d.Content\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $r.InsertAfter($text)\n}\n"}
